# Generate Report for Handoff
# Update the localization-status workbook so that the row for
# c617b04f-6424-4b11-adb3-db1762f1a537.md reflects that the file is now
# "Ready for handoff" (instead of "Handed back: in sync with en-US"),
# refresh the related timestamps, and record the "stale handback" error
# detail message on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfe93e339115c5a4322378e793c43696e8085dc6/e2e/c617b04f-6424-4b11-adb3-db1762f1a537.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b07dcc4025ff757673f0c5ba49add456b4c4d5d3/e2e/c617b04f-6424-4b11-adb3-db1762f1a537.md."

# ---------------------------------------------------------------------
# "Overview" sheet: row 3 corresponds to c617b04f-6424-4b11-adb3-db1762f1a537.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusReady
$overview.Range("F3").Value = $statusReady
$overview.Range("G3").Value = "2016-09-05 18:54:47"

# ---------------------------------------------------------------------
# "zh-cn" sheet: row 3 corresponds to c617b04f-6424-4b11-adb3-db1762f1a537.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReady
$zhcn.Range("H3").Value = "2016-09-05 18:54:41"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# "de-de" sheet: row 3 corresponds to c617b04f-6424-4b11-adb3-db1762f1a537.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReady
$dede.Range("H3").Value = "2016-09-05 18:54:47"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.17
